# "Brak" is the 4th sheet in the workbook (xl/worksheets/sheet4.xml)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Brak")

# New row of tracked time: 27.03.2018, "Use-Case Analyse & State Machine", 2h planned, 2h worked
# Copy the date cell's formatting (A2) down to A3 so it keeps the same date number format/style
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Cells.Item(3, 1).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(3, 1).Value = 43186
$ws.Cells.Item(3, 2).Value = "Use-Case Analyse & State Machine"
$ws.Cells.Item(3, 3).Value = 2
$ws.Cells.Item(3, 4).Value = 2

# Match row height of the other data rows
$ws.Rows.Item(3).RowHeight = 15.75

# Move the selection to B4, as in the saved file
$ws.Range("B4").Select() | Out-Null
